$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, reusing the existing header style (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for rows 2-14
$saveValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
